# Add two new columns "I0" (I) and "IF" (J) with header style matching
# the existing header row, then fill in the values for rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (copy style from an existing header cell, e.g. H1, so the new
# headers get the same bold/centered/bordered formatting)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-11
$i0Values = @(2, 1, 1, 1, 1, 1, 1, 7, 4, 1)
$ifValues = @(4, 4, 4, 5, 5, 5, 4, 9, 5, 2)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$r]
    $ws.Cells.Item($row, 10).Value = $ifValues[$r]
}

$wb.Save()
